$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the taxon-specific data between row 3 and row 4, while
# leaving the shared/location columns (C, I, P, S, T, U, V, W, Y, AA, AD,
# AE, AG, AS, AT, AW, AX, AY) untouched, since they are identical for
# both rows anyway.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"
    $v3 = $ws.Range($addr3).Value2
    $v4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value2 = $v4
    $ws.Range($addr4).Value2 = $v3
}
